$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new employee row (row 42): Pepe Pelele.
# Columns A-D are text, E/F are numbers, and none of the new cells carry
# any special style (unlike the existing rows, which center D/E/F).

# --- A42: must be stored as literal text "385580" (not a number) -------
# Writing the numeric-looking string directly would be auto-coerced to a
# number by the engine, so compute it via TEXT() and paste back as a
# value; this yields a plain inline/shared text cell with no number
# format applied (avoids leaving an unused style behind).
$ws.Cells.Item(42, 1).Formula = '=TEXT(385580,"0")'
$ws.Cells.Item(42, 1).Copy()
$ws.Cells.Item(42, 1).PasteSpecial(-4163)

# --- B42/C42/D42: plain text values -------------------------------------
$ws.Cells.Item(42, 2).Value = "Pepe"
$ws.Cells.Item(42, 3).Value = "Pelele"
$ws.Cells.Item(42, 4).Value = "A"

# --- E42/F42: plain numbers ---------------------------------------------
# Columns E/F default to a centered style for any new cell; reset back to
# Normal so row 42 matches (it has no style overrides).
$ws.Cells.Item(42, 5).Value = 1
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 6).Style = "Normal"
